$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","E","F","G","H","I","J","K","L","M")

$row8 = @("9 ماهه منتهی به 1399/09", "12 ماهه منتهی به 1399/12", "3 ماهه منتهی به 1400/03", "6 ماهه منتهی به 1400/06", "9 ماهه منتهی به 1400/09", "12 ماهه منتهی به 1400/12", "3 ماهه منتهی به 1401/03", "6 ماهه منتهی به 1401/06", "9 ماهه منتهی به 1401/09", "12 ماهه منتهی به 1401/12")
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "8").Value = $row8[$i] }

$row9 = @("1400-10-29 (2)", "1401-03-08 (8)", "1401-04-29 (2)", "1401-08-29 (4)", "1401-10-28 (2)", "1402-02-28 (7)", "1401-04-29", "1401-08-29 (2)", "1401-10-28", "1402-02-28")
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "9").Value = $row9[$i] }

$row11 = @(19427, 26191, 7048, 13918, 21257, 27409, 6369, 16281, 27368, 35012)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "11").Value = $row11[$i] }

$row12 = @(-8170, -11451, -3335, -7660, -11517, -15858, -3817, -7771, -13411, -17813)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "12").Value = $row12[$i] }

$row13 = @(11257, 14740, 3713, 6258, 9740, 11552, 2552, 8510, 13956, 17199)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "13").Value = $row13[$i] }

$row14 = @(-569, -856, -272, -598, -896, -1658, -426, -1100, -1193, -1653)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "14").Value = $row14[$i] }

$row15 = @("-", "-", "-", "-", "-", "-", "-", -78, -74, -202)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "15").Value = $row15[$i] }

$row16 = @(292, 296, -75, -25, -87, -98, 26, 105, -129, -121)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "16").Value = $row16[$i] }

$row17 = @(10981, 14180, 3366, 5634, 8757, 9796, 2151, 7438, 12561, 15224)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "17").Value = $row17[$i] }

$row18 = @(-1012, -1369, -423, -935, -1506, -2117, -670, -1537, -2364, -2867)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "18").Value = $row18[$i] }

$row19 = @(352, 352, 546, 556, 536, 535, 1036, 677, 320, 272)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "19").Value = $row19[$i] }

$row20 = @(10321, 13163, 3489, 5255, 7787, 8214, 2518, 6578, 10516, 12629)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "20").Value = $row20[$i] }

$row21 = @(-2322, -2154, -785, -1049, -1624, -1375, -333, -1225, -2366, -1861)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "21").Value = $row21[$i] }

$row22 = @(7999, 11009, 2704, 4206, 6163, 6839, 2184, 5353, 8150, 10768)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "22").Value = $row22[$i] }

$row23 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "23").Value = $row23[$i] }

$row24 = @(7999, 11009, 2704, 4206, 6163, 6839, 2184, 5353, 8150, 10768)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "24").Value = $row24[$i] }

$row25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "25").Value = $row25[$i] }

$row26 = @(4683, 6767, 5701, 6174, 5881, 9210, 8270, 8048, 11415, 10286)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "26").Value = $row26[$i] }

$row27 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i=0; $i -lt 10; $i++) { $ws.Range($cols[$i] + "27").Value = $row27[$i] }

# Column width updates (shift pattern for "wide" columns, offset 5/6 accounts for
# Excel's character-width <-> point-width conversion quirk so the stored XML width
# comes out as a clean integer)
$offset = 0.8333333333333334
$ws.Columns.Item(5).ColumnWidth = 29 - $offset
$ws.Columns.Item(6).ColumnWidth = 28 - $offset
$ws.Columns.Item(9).ColumnWidth = 29 - $offset
$ws.Columns.Item(10).ColumnWidth = 28 - $offset
$ws.Columns.Item(13).ColumnWidth = 29 - $offset
